$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D53").Value = 29991
$ws.Range("M53").Value = 43311.85572074554
$ws.Range("N53").Value = 59169.37404153174
$ws.Range("O53").Value = 16403.66459334077
$ws.Range("P53").Value = 29991
$ws.Range("Q53").Value = 1.132904101974189
$ws.Range("S53").Value = 0.07614414504696709
$ws.Range("T53").Value = 0.1040223588757164
$ws.Range("U53").Value = 0.02883836296810693
$ws.Range("V53").Value = 0.05272549550468167
$ws.Range("W53").Value = 1.031684161244126
$ws.Range("D101").Value = 17963.8
$ws.Range("E101").Value = 107749
$ws.Range("M101").Value = 17963.8
$ws.Range("N101").Value = 17963.8
$ws.Range("O101").Value = 5094.204749732103
$ws.Range("P101").Value = 6130.98976109215
$ws.Range("Q101").Value = 1.524500582130502
$ws.Range("R101").Value = 1.022228336147847
$ws.Range("S101").Value = 0.166718948667737
$ws.Range("T101").Value = 0.166718948667737
$ws.Range("U101").Value = 0.04727844109673503
$ws.Range("V101").Value = 0.05690066507431298
$ws.Range("W101").Value = 1.491350345340075
$ws.Range("D102").Value = 20152.1
$ws.Range("E102").Value = 41785
$ws.Range("M102").Value = 20152.1
$ws.Range("N102").Value = 20152.1
$ws.Range("O102").Value = 5714.766560364528
$ws.Range("P102").Value = 6877.849829351535
$ws.Range("Q102").Value = 1.289599257279878
$ws.Range("R102").Value = 0.8682958252810506
$ws.Range("S102").Value = 0.4822807227473974
$ws.Range("T102").Value = 0.4822807227473974
$ws.Range("U102").Value = 0.1367659820596991
$ws.Range("V102").Value = 0.1646009292653233
$ws.Range("W102").Value = 1.48520725279597
$ws.Range("D103").Value = 8250
$ws.Range("E103").Value = 56424
$ws.Range("M103").Value = 8250
$ws.Range("N103").Value = 8250
$ws.Range("O103").Value = 2339.548936488374
$ws.Range("P103").Value = 2815.699658703072
$ws.Range("Q103").Value = 0.7563944323516366
$ws.Range("R103").Value = 0.9503789792824658
$ws.Range("S103").Value = 0.1462143768609102
$ws.Range("T103").Value = 0.1462143768609102
$ws.Range("U103").Value = 0.04146371998597004
$ws.Range("V103").Value = 0.04990251769996937
$ws.Range("W103").Value = 0.7958871658995581
$ws.Range("D104").Value = 5608.1
$ws.Range("E104").Value = 30290
$ws.Range("M104").Value = 5608.1
$ws.Range("N104").Value = 5608.1
$ws.Range("O104").Value = 1590.354471602479
$ws.Range("P104").Value = 1914.027303754266
$ws.Range("Q104").Value = 0.7808335570377357
$ws.Range("R104").Value = 0.9381775382518739
$ws.Range("S104").Value = 0.1851469131726643
$ws.Range("T104").Value = 0.1851469131726643
$ws.Range("U104").Value = 0.05250427440087418
$ws.Range("V104").Value = 0.06319007275517552
$ws.Range("W104").Value = 0.8322876270227909
$ws.Range("D105").Value = 3155.9
$ws.Range("E105").Value = 23412
$ws.Range("M105").Value = 3155.9
$ws.Range("N105").Value = 3155.9
$ws.Range("O105").Value = 894.9554531713527
$ws.Range("P105").Value = 1077.098976109215
$ws.Range("Q105").Value = 0.6849351979224476
$ws.Range("R105").Value = 0.5617756448710258
$ws.Range("S105").Value = 0.1347983939859901
$ws.Range("T105").Value = 0.1347983939859901
$ws.Range("U105").Value = 0.03822635627760775
$ws.Range("V105").Value = 0.04600627781091812
$ws.Range("W105").Value = 1.21923263170246
$ws.Range("D106").Value = 24455.2
$ws.Range("M106").Value = 24455.2
$ws.Range("N106").Value = 24455.2
$ws.Range("O106").Value = 6935.046927467937
$ws.Range("P106").Value = 8346.484641638224
$ws.Range("Q106").Value = 2.587137002858179
$ws.Range("S106").Value = 1.335255255255255
$ws.Range("T106").Value = 1.335255255255255
$ws.Range("U106").Value = 0.3786539408936903
$ws.Range("V106").Value = 0.4557185171519642
$ws.Range("W106").Value = 1.6835107179942
$ws.Range("D107").Value = 48973.5
$ws.Range("E107").Value = 546163
$ws.Range("M107").Value = 48973.5
$ws.Range("N107").Value = 48973.5
$ws.Range("O107").Value = 13887.9878595289
$ws.Range("P107").Value = 16714.50511945392
$ws.Range("Q107").Value = 0.9376857527525851
$ws.Range("R107").Value = 1.054382974287297
$ws.Range("S107").Value = 0.08966828584140632
$ws.Range("T107").Value = 0.08966828584140632
$ws.Range("U107").Value = 0.02542828397296942
$ws.Range("V107").Value = 0.03060351052607724
$ws.Range("W107").Value = 0.8893217887802174
$ws.Range("D108").Value = 102721.4
$ws.Range("E108").Value = 536407
$ws.Range("M108").Value = 102721.4
$ws.Range("N108").Value = 102721.4
$ws.Range("O108").Value = 29129.90813631477
$ws.Range("P108").Value = 35058.49829351535
$ws.Range("Q108").Value = 1.279037367872705
$ws.Range("R108").Value = 0.9795938865665537
$ws.Range("S108").Value = 0.1914989923695999
$ws.Range("T108").Value = 0.1914989923695999
$ws.Range("U108").Value = 0.0543056077499264
$ws.Range("V108").Value = 0.06535801787358358
$ws.Range("W108").Value = 1.30568124751747
$ws.Range("D109").Value = 52429.3
$ws.Range("E109").Value = 230361
$ws.Range("M109").Value = 52429.3
$ws.Range("N109").Value = 52429.3
$ws.Range("O109").Value = 14867.98946131272
$ws.Range("P109").Value = 17893.9590443686
$ws.Range("Q109").Value = 1.331258261821052
$ws.Range("R109").Value = 0.8644815797472915
$ws.Range("S109").Value = 0.2275962511015319
$ws.Range("T109").Value = 0.2275962511015319
$ws.Range("U109").Value = 0.06454212935919151
$ws.Range("V109").Value = 0.07767790139984025
$ws.Range("W109").Value = 1.539949830059086
$ws.Range("D143").Value = 112450
$ws.Range("M143").Value = 61880.92159065686
$ws.Range("N143").Value = 70953.65163674846
$ws.Range("O143").Value = 20195.69439758588
$ws.Range("P143").Value = 35326.45688255929
$ws.Range("Q143").Value = 1.358535294507217
$ws.Range("S143").Value = 0.08726881665734507
$ws.Range("T143").Value = 0.1000638170326061
$ws.Range("U143").Value = 0.02848138499470568
$ws.Range("V143").Value = 0.04981984769443295
$ws.Range("W143").Value = 0.9924237083443946
$ws.Range("D161").Value = 174376
$ws.Range("M161").Value = 68034.78193023261
$ws.Range("N161").Value = 81722.72026718804
$ws.Range("O161").Value = 22463.64548950225
$ws.Range("P161").Value = 44826.735218509
$ws.Range("Q161").Value = 1.564728485216028
$ws.Range("S161").Value = 0.08673038336932412
$ws.Range("T161").Value = 0.1041796954097038
$ws.Range("U161").Value = 0.0286365375165752
$ws.Range("V161").Value = 0.05714488707677961
$ws.Range("W161").Value = 1.033244610476908
